$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Error [mg/l]" and "Conc [mg/l]" columns (B and C) were swapped,
# including their headers and all data values, for rows 1 through 5.
for ($r = 1; $r -le 5; $r++) {
    $bCell = $ws.Cells.Item($r, 2)
    $cCell = $ws.Cells.Item($r, 3)
    $bVal = $bCell.Value()
    $cVal = $cCell.Value()
    $bCell.Value = $cVal
    $cCell.Value = $bVal
}
